$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14); this pushes the
# existing "Late" / "heading" / "Outstanding" columns one to the right
# (N->O, O->P, P->Q), matching the target layout.
$ws.Columns.Item(14).Insert()

# The newly inserted column picks up the default sheet width; give it
# the same width as the column to its left (M / column 13), which is
# what Excel does when inserting a column in the middle of a table.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab and move the
# selection to L18 on it (previously the Transactions sheet was
# active with selection C3, which stays untouched on that sheet).
$ws.Range("L18").Select()
$ws.Activate()
